$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before current row 3 (Spain), shifting Spain/Italy down to rows 7/8
$ws.Range("A3:A6").EntireRow.Insert()

function Set-Row($r, $values) {
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $values[$i]
    }
}

# Row 2 - Germany, week 22
Set-Row 2 @("Germany", 22, 0.06464924231864996, 0.001925425130163749, 0.0004161797637778281, 0.001509245366385914, 0.216149543941205, 0.783850456058795)

# Row 3 - Germany, week 21
Set-Row 3 @("Germany", 21, 0.06513333139719893, 0.001441336051614783, 0.0001703283534777427, 0.001271007698137038, 0.1181739354170165, 0.8818260645829835)

# Row 4 - Germany, week 20
Set-Row 4 @("Germany", 20, 0.06532623255231347, 0.001248434896500242, 0.00006138049671319686, 0.001187054399787047, 0.04916595721992851, 0.9508340427800716)

# Row 5 - USA, week 21
Set-Row 5 @("USA", 21, 0.06563651825795751, 0.0009381491908562051, 0.004017276000219924, -0.003079126809363722, 0.5661003339318083, 0.4338996660681918)

# Row 6 - Germany, week 19 (only A-F populated, G/H left blank)
Set-Row 6 @("Germany", 19, 0.06657466744881371, 0, 0, 0)

# Row 7 - Spain, week 19 (updated D/E/F/G/H values)
Set-Row 7 @("Spain", 19, 0.2179061588686753, -0.1513314914198616, -0.07677906293842739, -0.07455242848143416, 0.5073568113156817, 0.4926431886843184)

# Row 8 - Italy, week 19 (updated D/E/F/G/H values)
Set-Row 8 @("Italy", 19, 0.2250579911111283, -0.1584833236623146, -0.08887279934180774, -0.06961052432050688, 0.5607706684090738, 0.4392293315909263)
